$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): columns A..W
$headers = @(
    "productIds",
    "MSE_no_transfer",
    "MSE_transfer_basic",
    "MSE_transfer_coral",
    "MSE_transfer_sa",
    "MSE_transfer_bw",
    "MSE_transfer_nnw",
    "MAE_no_transfer",
    "MAE_transfer_basic",
    "MAE_transfer_coral",
    "MAE_transfer_sa",
    "MAE_transfer_bw",
    "MAE_transfer_nnw",
    "MSE_diff_basic",
    "MSE_transfer_coral",
    "MSE_diff_sa",
    "MSE_diff_bw",
    "MSE_diff_nnw",
    "MAE_diff_basic",
    "MAE_transfer_coral",
    "MAE_diff_sa",
    "MAE_diff_bw",
    "MAE_diff_nnw"
)

for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Data rows 2..6, columns B..W (column A / productIds unchanged)
$data = @{
    2 = @(1.224559961485893, 1.227225739489796, 1.120835828974163, 1.033296906461127, 0.8809340079803796, 0.890879775520842, 0.5149422876031311, 0.5028908839856228, 0.5656839662244296, 0.78726599205821, 0.4370141014276453, 0.3903746227792839, 0.002665778003903441, -0.10372413251173, -0.191263055024766, -0.3436259535055131, -0.3336801859650507, -0.01205140361750834, 0.05074167862129852, 0.2723237044550789, -0.0779281861754858, -0.1245676648238472)
    3 = @(2.588536900271843, 3.013620179988538, 3.161065882811345, 2.852513077306079, 2.806361749175249, 2.699051178012179, 0.9914813140809273, 0.8658876819259795, 0.8625860071605033, 1.032059815070777, 0.8616075054061442, 0.7461395487668154, 0.4250832797166946, 0.5725289825395015, 0.2639761770342361, 0.2178248489034056, 0.1105142777403354, -0.1255936321549478, -0.1288953069204241, 0.04057850098984983, -0.1298738086747832, -0.2453417653141119)
    4 = @(5.485710658411413, 5.71663333954605, 3.502699856861465, 3.176874955584172, 3.230494342152407, 3.11746815953909, 1.282621473071498, 1.290585842418499, 1.009296810350896, 1.108763838657743, 0.9005700749771435, 0.8744593543504965, 0.2309226811346372, -1.983010801549948, -2.30883570282724, -2.255216316259006, -2.368242498872323, 0.007964369347001066, -0.2733246627206014, -0.1738576344137543, -0.3820513980943543, -0.4081621187210013)
    5 = @(14.22789284129596, 13.62535385482689, 12.53143968824912, 12.63280682930526, 12.37904794826922, 12.56539049160015, 1.136556738826326, 1.075380537244296, 1.015973113539503, 1.018982735495265, 0.9943578763063982, 0.9724489795748279, -0.6025389864690691, -1.696453153046841, -1.595086011990695, -1.848844893026735, -1.662502349695812, -0.06117620158202919, -0.1205836252868224, -0.1175740033310602, -0.1421988625199274, -0.1641077592514977)
    6 = @(8.201373654786686, 7.811590549391124, 5.859342160701627, 5.807247552767949, 5.88428453792401, 6.225046535432612, 1.062188837815168, 1.013146724023256, 0.8617176915191481, 0.8229994240953901, 0.808653786169095, 0.7944319410373135, -0.3897831053955612, -2.342031494085059, -2.394126102018737, -2.317089116862675, -1.976327119354074, -0.04904211379191192, -0.2004711462960199, -0.239189413719778, -0.253535051646073, -0.2677568967778545)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($r, $c + 2).Value = $vals[$c]
    }
}

# The new header cells beyond the original used range (L1:W1) need the
# same header formatting (bold, centered, thin border) as the existing
# header cells. Copy the format from an existing styled header cell.
$ws.Range("A1").Copy()
$ws.Range("L1:W1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

Write-Output "applied"
